{"js": "// The \"Discussion\" paragraph mentions the winning entrant's approach:\n//   \"...was able to reach an overall test accuracy of over 90% using\n//    ArcFace (Additive Angular Margin Loss for Deep Face Recognition).\n//    The approach of using convolutional neural nets...\"\n//\n// The commit adds another finalist model to that list, turning the\n// sentence boundary \"Recognition). The approach\" into\n// \"Recognition) and senet154. The approach\", i.e. it inserts\n// \" and senet154\" right after the closing parenthesis of the ArcFace\n// citation and before the full stop that starts the next sentence.\n\nconst body = context.document.body;\n\n// Narrow, unambiguous anchor: the single place in the document where the\n// ArcFace citation's closing paren is immediately followed by the next\n// sentence.\nconst searchText = \"Recognition). The approach of using\";\nconst results = body.search(searchText, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Could not find target text \"${searchText}\" in document body.`);\n}\n\n// Anchor on just the \")\" so we can insert immediately after it, leaving\n// \". The approach of using\" untouched.\nconst parenResults = body.search(\"Recognition)\", { matchCase: true });\nparenResults.load(\"text\");\nawait context.sync();\n\nif (parenResults.items.length === 0) {\n  throw new Error('Could not find anchor text \"Recognition)\" in document body.');\n}\n\nconst anchor = parenResults.items[0];\nanchor.insertText(\" and senet154\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# The \"Discussion\" paragraph mentions the winning entrant's approach:\n#   \"...was able to reach an overall test accuracy of over 90% using\n#    ArcFace (Additive Angular Margin Loss for Deep Face Recognition).\n#    The approach of using convolutional neural nets...\"\n#\n# The commit adds another finalist model to that list, turning the\n# sentence boundary \"Recognition). The approach\" into\n# \"Recognition) and senet154. The approach\", i.e. it inserts\n# \" and senet154\" right after the closing parenthesis of the ArcFace\n# citation and before the full stop that starts the next sentence.\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.MatchCase = $true\n$range.Find.MatchWildcards = $false\n$range.Find.Execute(\"Recognition)\") | Out-Null\n\nif (-not $range.Find.Found) {\n    throw \"Could not find anchor text 'Recognition)' in document body.\"\n}\n\n# $range now covers exactly \"Recognition)\" - collapse to its end point and\n# insert the new text right after it, leaving the following\n# \". The approach of using \" text untouched.\n$range.Collapse(0)  # wdCollapseEnd\n$range.InsertAfter(\" and senet154\")\n"}
